$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section 1/2: Indicator + Organization contact details were refreshed
# (new responsible person at the National Statistical Committee).

# Индикатор (same text, now with a trailing space as published)
$ws.Range("B4").Value = "5.3.1. Доля женщин в возрасте от 20 до 24 лет, вступивших в брак или союз до 15 лет и до 18 лет "

# Сайт организации (если есть)
$ws.Range("B10").Value = "www.stat.gov.kg"

# Телефон контактного лица
$ws.Range("B9").Value = "(0312) 32 46 55"

# Электронная почта контактного лица
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com"

# Организация: "Отдел" -> "Управление"
$ws.Range("B6").Value = "Национальный статистический комитет КР" + [char]10 + "(Управление статистики домашних хозяйств)"

# Контактное лицо (лица) / Координатор
$ws.Range("B7").Value = "Калымбетова Ы.И."

# Match the author's last on-screen selection when the file was saved.
$ws.Range("B8").Select()
